$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Phase 1: copy cell FORMATS from existing template rows onto the new row
# positions (rows 89-100) while the template rows are still untouched, and
# fix up the few format changes needed inside the existing table block
# (rows 83, 84, 87) using other untouched template rows. Values are set
# afterwards in Phase 2, so it does not matter that some of these templates
# get overwritten with new values later.
# ---------------------------------------------------------------------------

$ws.Range("A72:E72").Copy()
$ws.Range("A100:E100").PasteSpecial(-4122)

$ws.Range("A74:D74").Copy()
$ws.Range("A89:D89").PasteSpecial(-4122)

$ws.Range("A75:D75").Copy()
$ws.Range("A90:D90").PasteSpecial(-4122)

$ws.Range("A76:E76").Copy()
$ws.Range("A91:E91").PasteSpecial(-4122)

$ws.Range("A77:E77").Copy()
$ws.Range("A92:E92").PasteSpecial(-4122)

$ws.Range("A78:E78").Copy()
$ws.Range("A93:E93").PasteSpecial(-4122)

$ws.Range("A83:E83").Copy()
$ws.Range("A94:E94").PasteSpecial(-4122)

$ws.Range("A84:E84").Copy()
$ws.Range("A95:E95").PasteSpecial(-4122)

$ws.Range("A85:E85").Copy()
$ws.Range("A96:E96").PasteSpecial(-4122)

$ws.Range("A86:E86").Copy()
$ws.Range("A97:E97").PasteSpecial(-4122)

$ws.Range("A87:E87").Copy()
$ws.Range("A98:E98").PasteSpecial(-4122)

$ws.Range("A88:E88").Copy()
$ws.Range("A99:E99").PasteSpecial(-4122)

# Fix formatting for the rows that stay inside the T_MATCH_FOOTBALL_MEMBER
# block but change column D's style (empty bordered cell <-> shared-string
# bordered cell).
$ws.Range("A85:D85").Copy()
$ws.Range("A83:D83").PasteSpecial(-4122)

$ws.Range("A86:D86").Copy()
$ws.Range("A84:D84").PasteSpecial(-4122)

$ws.Range("A72:D72").Copy()
$ws.Range("A87:D87").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Phase 2: write the actual values / formulas for rows 82-100
# ---------------------------------------------------------------------------

# Row 82 : MATCH_MEMBER_POSITION VARCHAR(100) ,
$ws.Range("B82").Value = "VARCHAR(100)"

# Row 83 : INPUT_NAME VARCHAR(100) ,
$ws.Range("A83").Value = "INPUT_NAME"
$ws.Range("B83").Value = "VARCHAR(100)"
$ws.Range("D83").Value = ""

# Row 84 : INPUT_DATE DATE ,
$ws.Range("A84").Value = "INPUT_DATE"
$ws.Range("B84").Value = "DATE"
$ws.Range("D84").Value = ""

# Row 85 : UPDATE_NAME VARCHAR(100) ,
$ws.Range("A85").Value = "UPDATE_NAME"

# Row 86 : UPDATE_DATE DATE ,
$ws.Range("A86").Value = "UPDATE_DATE"

# Row 87 : DELSIGN DECIMAL(1,0) ,
$ws.Range("A87").Value = "DELSIGN"
$ws.Range("B87").Value = "DECIMAL(1,0)"
$ws.Range("D87").Value = "0:UNDELETED;1:DELETED;"

# Row 89 : NO / DATE header for the new T_DICTIONARY table block
$ws.Range("A89").Value = "NO"
$ws.Range("B89").Value = 5
$ws.Range("C89").Value = "DATE"
$ws.Range("D89").Value = 42269

# Row 90 : TABLE_NAME = T_DICTIONARY
$ws.Range("A90").Value = "TABLE_NAME"
$ws.Range("B90").Value = "T_DICTIONARY"
$ws.Range("C90").Value = ""
$ws.Range("D90").Value = ""

# Row 91 : AUTHOR = Rex , drop table formula
$ws.Range("A91").Value = "AUTHOR"
$ws.Range("B91").Value = "Rex"
$ws.Range("C91").Value = ""
$ws.Range("D91").Value = ""
$ws.Range("E91").Formula = '="drop table "&B90&";"'

# Row 92 : COLUMN / TYPE / CONSTRAINTS / COMMENT header , create table formula
$ws.Range("A92").Value = "COLUMN"
$ws.Range("B92").Value = "TYPE"
$ws.Range("C92").Value = "CONSTRAINTS"
$ws.Range("D92").Value = "COMMENT"
$ws.Range("E92").Formula = '="create table "&B90&" ("'

# Row 93 : D_ID DECIMAL(10,0) NOT NULL PRIMARY KEY,
$ws.Range("A93").Value = "D_ID"
$ws.Range("B93").Value = "DECIMAL(10,0)"
$ws.Range("C93").Value = "NOT NULL PRIMARY KEY"
$ws.Range("D93").Value = ""
$ws.Range("E93").Formula = '=A93&" "&B93&" "&C93&", "'

# Row 94 : D_NAMME VARCHAR(500) , 名称
$ws.Range("A94").Value = "D_NAMME"
$ws.Range("B94").Value = "VARCHAR(500)"
$ws.Range("C94").Value = ""
$ws.Range("D94").Value = "名称"
$ws.Range("E94").Formula = '=A94&" "&B94&" "&C94&", "'

# Row 95 : D_TYPE DECIMAL(10,0) , 类型
$ws.Range("A95").Value = "D_TYPE"
$ws.Range("B95").Value = "DECIMAL(10,0)"
$ws.Range("C95").Value = ""
$ws.Range("D95").Value = "类型"
$ws.Range("E95").Formula = '=A95&" "&B95&" "&C95&", "'

# Row 96 : INPUT_NAME VARCHAR(100) ,
$ws.Range("A96").Value = "INPUT_NAME"
$ws.Range("B96").Value = "VARCHAR(100)"
$ws.Range("C96").Value = ""
$ws.Range("D96").Value = ""
$ws.Range("E96").Formula = '=A96&" "&B96&" "&C96&", "'

# Row 97 : INPUT_DATE DATE ,
$ws.Range("A97").Value = "INPUT_DATE"
$ws.Range("B97").Value = "DATE"
$ws.Range("C97").Value = ""
$ws.Range("D97").Value = ""
$ws.Range("E97").Formula = '=A97&" "&B97&" "&C97&", "'

# Row 98 : UPDATE_NAME VARCHAR(100) ,
$ws.Range("A98").Value = "UPDATE_NAME"
$ws.Range("B98").Value = "VARCHAR(100)"
$ws.Range("C98").Value = ""
$ws.Range("D98").Value = ""
$ws.Range("E98").Formula = '=A98&" "&B98&" "&C98&", "'

# Row 99 : UPDATE_DATE DATE ,
$ws.Range("A99").Value = "UPDATE_DATE"
$ws.Range("B99").Value = "DATE"
$ws.Range("C99").Value = ""
$ws.Range("D99").Value = ""
$ws.Range("E99").Formula = '=A99&" "&B99&" "&C99&", "'

# Row 100 : DELSIGN DECIMAL(1,0) ,
$ws.Range("A100").Value = "DELSIGN"
$ws.Range("B100").Value = "DECIMAL(1,0)"
$ws.Range("C100").Value = ""
$ws.Range("D100").Value = "0:UNDELETED;1:DELETED;"
$ws.Range("E100").Formula = '=A100&" "&B100&" "&C100&", "'

# ---------------------------------------------------------------------------
# Phase 3: sheet view selection to match the edited workbook
# ---------------------------------------------------------------------------
$ws.Range("C103").Select()
